# Remove the visible attendance "password" placeholder from the slide and
# replace it with "moreFun!" (so the instructor tells students to look on
# the board for the real password instead of showing it on the slide).

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null
$targetParaIndex = 0

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -like "*___________*") {
                    $targetSlide = $s
                    $targetShape = $shp
                    # find which paragraph holds the blank-password placeholder
                    $paraCount = $tr.Paragraphs().Count
                    for ($pi = 1; $pi -le $paraCount; $pi++) {
                        $para = $tr.Paragraphs($pi, 1)
                        if ($para.Text -like "*___________*") {
                            $targetParaIndex = $pi
                        }
                    }
                }
            }
        }
    }
}

if ($targetShape -ne $null -and $targetParaIndex -gt 0) {
    $tr = $targetShape.TextFrame.TextRange
    $para = $tr.Paragraphs($targetParaIndex, 1)

    # Replace the blank "___________" with "moreFun!" while keeping the
    # paragraph's existing run formatting (44pt, yellow highlight, etc.)
    $para.Text = "moreFun!"

    # Re-assert the text of the first 7 characters ("moreFun") so that it
    # becomes its own run, split off from the trailing "!" -- mirroring how
    # the two segments ended up as separate runs in the authored slide.
    $firstPart = $para.Characters(1, 7)
    $firstPart.Text = "moreFun"
}
